$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = ' NOROESTE DE LAS ISLAS KURILES '
$ws.Range("A3").Value = ' AL SUR DE LAS ISLAS FIJI '
$ws.Range("A4").Value = ' ALASKA PENINSULA '
$ws.Range("A5").Value = ' TIMOR REGIÓN '
$ws.Range("A6").Value = ' JÓNICO '
$ws.Range("A7").Value = ' NORTE DE ARGELIA '
$ws.Range("A8").Value = ' OESTE CHILE RISE '
$ws.Range("A9").Value = ' CHILE-BOLIVIA FRONTERA REGIÓN '
$ws.Range("A10").Value = ' HOKKAIDO, JAPAN REGION '
$ws.Range("A11").Value = ' AL SUR DE LAS ISLAS FIJI '
$ws.Range("A12").Value = ' NUEVA GUINEA, PAPUA NUEVA GUINEA '
$ws.Range("A13").Value = ' CANAL DE LA MONA '
$ws.Range("A14").Value = ' CENTRAL DE ITALIA '
$ws.Range("A15").Value = ' FLORES DEL MAR '
$ws.Range("A16").Value = ' REGIÓN DE LAS ISLAS SAMOA '
$ws.Range("A17").Value = ' KAZAJSTáN ORIENTAL '
$ws.Range("A18").Value = ' CERCA DE COSTA DE NICARAGUA '
$ws.Range("A19").Value = ' AFGANISTÁN-TAYIKISTÁN FRONT. REGIÓN '
$ws.Range("A20").Value = ' TAIWÁN '
$ws.Range("A21").Value = ' CERCA LA COSTA ESTE DE KAMCHATKA '
$ws.Range("A22").Value = ' ISLA NORTE, NUEVA ZELANDA '
$ws.Range("A23").Value = ' ISLAS SUR DE SANDWICH REGIÓN '
$ws.Range("A24").Value = ' MAR BEAUFORT '
$ws.Range("A25").Value = ' ISLAS KURILES '
$ws.Range("A26").Value = ' JAVA, INDONESIA '
$ws.Range("A27").Value = ' ISLAS SUR DE SANDWICH REGIÓN '
$ws.Range("A28").Value = ' CHILE-BOLIVIA FRONTERA REGIÓN '
$ws.Range("A29").Value = ' CERCA LA COSTA ESTE DE KAMCHATKA '
$ws.Range("A30").Value = ' SICHUAN, CHINA '
$ws.Range("A31").Value = ' BAFFIN REGIÓN DE LA ISLA, CANADA '
$ws.Range("A32").Value = ' ISLAS FIJI REGIÓN '
$ws.Range("A33").Value = ' ISLA DE ANDREANOF, ISLAS ALEUTIANAS '
$ws.Range("A34").Value = ' ISLAS KURILES '
$ws.Range("A35").Value = ' AFGANISTÁN-TAYIKISTÁN FRONT. REGIÓN '
$ws.Range("A36").Value = ' CERCA LA COSTA CENTRAL DE CHILE '
$ws.Range("A37").Value = ' CHIPRE REGIóN '
$ws.Range("A38").Value = ' ESTE DEL LAGO BAIKAL, RUSIA '
$ws.Range("A39").Value = ' CERCA LA COSTA DE PERÚ '
$ws.Range("A40").Value = ' PANAMÁ-COSTA RICA REGIÓN FRONTERA '
$ws.Range("A41").Value = ' CERCA LA COSTA CENTRAL DE CHILE '
$ws.Range("A42").Value = ' CERCA LA COSTA CENTRAL DE CHILE '
$ws.Range("A43").Value = ' REGIÓN DE LAS ISLAS KERMADAC '
$ws.Range("A44").Value = ' LUZON, FILIPINAS '
$ws.Range("A45").Value = ' ESTE MAR DE JAPÓN '
$ws.Range("A46").Value = ' SUROESTE DE ÁFRICA '
$ws.Range("A47").Value = ' ESTRECHO DE GIBRALTAR '
$ws.Range("A48").Value = ' EL SURESTE DE AFGANISTáN '
$ws.Range("A49").Value = ' YUNNAN, CHINA '
$ws.Range("A50").Value = ' CRETA, GRECIA '
$ws.Range("A51").Value = ' CERCA DE LA COSTA DE OREGON '
$ws.Range("A52").Value = ' SOUTHWEST OF SUMATRA, INDONESIA '
$ws.Range("A53").Value = ' SOUTHWEST OF SUMATRA, INDONESIA '
$ws.Range("A54").Value = ' SICHUAN, CHINA '
$ws.Range("A55").Value = ' ISLAS VANUATU '
$ws.Range("A56").Value = ' SUR DE XINJIANG, CHINA '
$ws.Range("A57").Value = ' AL SURESTE DE LAS ISLAS DE LEALTAD (LOYALTY) '
$ws.Range("A58").Value = ' MAR BANDA '
$ws.Range("A59").Value = ' CERCA LA COSTA ESTE DE KAMCHATKA '
$ws.Range("A60").Value = ' ISLAS KERMADAC, NUEVA ZELANDA '
$ws.Range("A61").Value = ' ESTRECHO DE SUNDA, INDONESIA '
$ws.Range("A62").Value = ' ISLAS SANTA CRUZ '
$ws.Range("A63").Value = ' AL SUR DE LAS ISLAS KERMADAC '
$ws.Range("A64").Value = ' TAYIKISTÁN-XINJIANG REGIÓN FRONTERA '
$ws.Range("A65").Value = ' SUR DE ALASKA '
$ws.Range("A66").Value = ' ALASKA CENTRAL '
$ws.Range("A67").Value = ' HALMAHERA, INDONESIA '
$ws.Range("A68").Value = ' EN EL NORTE DE CHILE '
$ws.Range("A69").Value = ' CERCA DE COSTA DE GUATEMALA '
$ws.Range("A70").Value = ' CERCA LA COSTA E. DE HONSHU, JAPóN '
$ws.Range("A71").Value = ' N DE LAS ISLAS ARU REGIÓ INDONESIA '
$ws.Range("A72").Value = ' EL SUR DE IRáN '
$ws.Range("A73").Value = ' CERCA LA COSTA E. DE HONSHU, JAPóN '
$ws.Range("A74").Value = ' ISLAS FIJI REGIÓN '
$ws.Range("A75").Value = ' CANAL DE LA MONA '
$ws.Range("A76").Value = ' AL SURESTE DE LAS ISLAS DE LEALTAD (LOYALTY) '
$ws.Range("A77").Value = ' AL SUR DE LAS ISLAS FIJI '
$ws.Range("A78").Value = ' PAKISTÁN '
$ws.Range("A79").Value = ' NORTE DE MAR DE LAS MOLUCAS '
$ws.Range("A80").Value = ' REGIÓN NUEVA IRLANDA, PNG '
$ws.Range("A81").Value = ' ISLAS SALOMÓN '
$ws.Range("A82").Value = ' ECUADOR '
$ws.Range("A83").Value = ' EL NORTE DE XINJIANG, CHINA '
$ws.Range("A84").Value = ' EN EL NORTE DE CHILE '
$ws.Range("A85").Value = ' CERCA LA COSTA E. DE HONSHU, JAPóN '
$ws.Range("A86").Value = ' ISLAS SANTA CRUZ '
$ws.Range("A87").Value = ' EL OESTE DE TEXAS '
$ws.Range("A88").Value = ' EN EL NORTE DE CHILE '
$ws.Range("A89").Value = ' EN EL NORTE DE INDIA '
$ws.Range("A90").Value = ' REGIóN DE TAIWáN '
$ws.Range("A91").Value = ' REGIóN DE TAIWáN '
$ws.Range("A92").Value = ' ISLAS FIJI REGIÓN '
$ws.Range("A93").Value = ' OAXACA, MÉXICO '
$ws.Range("A94").Value = ' ISLAS FIJI REGIÓN '
$ws.Range("A95").Value = ' SUR DE PANAMÁ '
$ws.Range("A96").Value = ' SOUTHWEST OF SUMATRA, INDONESIA '
$ws.Range("A97").Value = ' E. RUSIA-NE CHINA FRONTERA REGIÓN '
$ws.Range("A98").Value = ' SICHUAN, CHINA '
$ws.Range("A99").Value = ' TAIWÁN '
$ws.Range("A100").Value = ' VOLCAN ISLAS, JAPÓN REGIÓN '
$ws.Range("A101").Value = ' ISLAS SUR DE SANDWICH REGIÓN '
